$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column stays text (avoid Excel auto-converting numeric-looking
# strings like "254.20" or "0.608" into floating point numbers, which would
# lose trailing zeros / introduce binary rounding artifacts).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.436.57"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "2.187.86"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "254.20"
$ws.Range("E5").Value = "  +4.50%  "
$ws.Range("D6").Value = "0.608"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("D7").Value = "74.26"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").Value = "  -4.38%  "
$ws.Range("D10").Value = "40.34"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "6.78"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "2.519.15"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "14.26"
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("D16").Value = "2.189.05"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "0.773"
$ws.Range("E17").Value = "  -4.77%  "
$ws.Range("D18").Value = "42.358.14"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "70.96"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "5.87"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "227.09"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "9.42"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "10.50"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").Value = "3.34"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "171.81"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "36.94"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "20.07"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").Value = "0.0826"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").Value = "5.14"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0340"
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "4.20"
$ws.Range("E38").Value = "  -5.25%  "
$ws.Range("D39").Value = "11.99"
$ws.Range("E39").Value = "  -7.22%  "
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").Value = "0.196"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").Value = "59.22"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  +10.89%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "5.14"
$ws.Range("E44").Value = "  -7.07%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "101.72"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").Value = "0.0971"
$ws.Range("D47").Value = "0.461"
$ws.Range("E47").Value = "  +7.07%  "
$ws.Range("D48").Value = "8.19"
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("D49").Value = "1.09"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "2.64"
$ws.Range("E51").Value = "  +0.24%  "
